$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 101
$ws.Range("I8").Value = 200
$ws.Range("K8").Value = 600
$ws.Range("M8").Value = -461
$ws.Range("H33").Value = 292.11765
$ws.Range("I33").Value = 223.27272
$ws.Range("K33").Value = 223.27272
$ws.Range("M33").Value = 5.727280000000007
$ws.Range("H92").Value = 440.1111
$ws.Range("I92").Value = 394.5
$ws.Range("K92").Value = 394.5
$ws.Range("M92").Value = 853.5
$ws.Range("H112").Value = 2222.0417
$ws.Range("J112").Value = 2425.3333
$ws.Range("L112").Value = 7275.999899999999
$ws.Range("N112").Value = -9491.999899999999
$ws.Range("H113").Value = 3195
$ws.Range("J113").Value = 3195
$ws.Range("L113").Value = 3195
$ws.Range("N113").Value = -9703
$ws.Range("H138").Value = 4178.857
$ws.Range("J138").Value = 3542.25
$ws.Range("L138").Value = 10626.75
$ws.Range("N138").Value = -20906.75
$ws.Range("H141").Value = 6624.875
$ws.Range("I141").Value = 6624.875
$ws.Range("K141").Value = 19874.625
$ws.Range("M141").Value = -14694.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 229.08333
$ws.Range("I2").Value = 229.08333
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 229.08333
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -116.08333
$ws.Range("N2").ClearContents()
$ws.Range("H45").Value = 6999.75
$ws.Range("I45").Value = 6999
$ws.Range("K45").Value = 6999
$ws.Range("M45").Value = -6622
$ws.Range("H61").Value = 7290.125
$ws.Range("I61").Value = 7290.125
$ws.Range("K61").Value = 7290.125
$ws.Range("M61").Value = -7078.125
$ws.Range("H116").Value = 229.08333
$ws.Range("I116").Value = 229.08333
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 229.08333
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 2064.91667
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 2023.8334
$ws.Range("I122").Value = 1985.75
$ws.Range("J122").Value = 2100
$ws.Range("K122").Value = 5957.25
$ws.Range("L122").Value = 6300
$ws.Range("M122").Value = -3507.25
$ws.Range("N122").Value = -11200
$ws.Range("H132").Value = 1543.4706
$ws.Range("I132").Value = 1387.6428
$ws.Range("J132").Value = 2270.6667
$ws.Range("K132").Value = 4162.928400000001
$ws.Range("L132").Value = 6812.000100000001
$ws.Range("M132").Value = -1632.928400000001
$ws.Range("N132").Value = -11872.0001
$ws.Range("H134").Value = 80000
$ws.Range("J134").Value = 80000
$ws.Range("L134").Value = 80000
$ws.Range("N134").Value = -90140
$ws.Range("H136").Value = 7290.125
$ws.Range("I136").Value = 7290.125
$ws.Range("K136").Value = 21870.375
$ws.Range("M136").Value = -19320.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 229.08333
$ws.Range("I3").Value = 229.08333
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 229.08333
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -115.08333
$ws.Range("N3").ClearContents()
$ws.Range("H22").Value = 645.25
$ws.Range("I22").Value = 753.6667
$ws.Range("J22").Value = 320
$ws.Range("K22").Value = 753.6667
$ws.Range("L22").Value = 320
$ws.Range("M22").Value = -580.6667
$ws.Range("N22").Value = -666
$ws.Range("H64").Value = 1075
$ws.Range("J64").Value = 1075
$ws.Range("L64").Value = 1075
$ws.Range("N64").Value = -1525
$ws.Range("H67").Value = 1075
$ws.Range("J67").Value = 1075
$ws.Range("L67").Value = 1075
$ws.Range("N67").Value = -2635

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4320.9287
$ws.Range("I31").Value = 3043.7273
$ws.Range("K31").Value = 3043.7273
$ws.Range("M31").Value = -2748.7273
$ws.Range("H34").Value = 4320.9287
$ws.Range("I34").Value = 3043.7273
$ws.Range("K34").Value = 3043.7273
$ws.Range("M34").Value = -2841.7273
$ws.Range("H41").Value = 12585.286
$ws.Range("I41").Value = 3619.6
$ws.Range("K41").Value = 3619.6
$ws.Range("M41").Value = -3191.6
$ws.Range("H122").Value = 910.1429000000001
$ws.Range("I122").Value = 910.1429000000001
$ws.Range("K122").Value = 2730.4287
$ws.Range("M122").Value = -280.4287000000004
$ws.Range("H132").Value = 6742.7
$ws.Range("I132").Value = 6103
$ws.Range("K132").Value = 18309
$ws.Range("M132").Value = -15779

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 83199.39999999999
$ws.Range("J22").Value = 55000
$ws.Range("L22").Value = 165000
$ws.Range("N22").Value = -165338
$ws.Range("H27").Value = 83199.39999999999
$ws.Range("J27").Value = 55000
$ws.Range("L27").Value = 165000
$ws.Range("N27").Value = -165204
$ws.Range("H34").Value = 617.4286
$ws.Range("J34").Value = 1033
$ws.Range("L34").Value = 3099
$ws.Range("N34").Value = -3267
$ws.Range("H60").Value = 322.5
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H92").Value = 999.875
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 999.875
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 2999.625
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -5495.625
$ws.Range("H104").Value = 9333.333000000001
$ws.Range("J104").Value = 9333.333000000001
$ws.Range("L104").Value = 27999.999
$ws.Range("N104").Value = -33241.999
$ws.Range("H121").Value = 2986
$ws.Range("I121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("H122").Value = 735.4
$ws.Range("J122").Value = 992.3333
$ws.Range("L122").Value = 8930.9997
$ws.Range("N122").Value = -13830.9997
$ws.Range("H131").Value = 922.9231
$ws.Range("J131").Value = 999.8
$ws.Range("L131").Value = 2999.4
$ws.Range("N131").Value = -13079.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 75995
$ws.Range("J93").Value = 75995
$ws.Range("L93").Value = 75995
$ws.Range("N93").Value = -79739
$ws.Range("H97").Value = 973.9091
$ws.Range("I97").Value = 982.7143
$ws.Range("J97").Value = 958.5
$ws.Range("K97").Value = 982.7143
$ws.Range("L97").Value = 958.5
$ws.Range("M97").Value = -486.7143
$ws.Range("N97").Value = -1950.5
$ws.Range("H102").Value = 1864.9286
$ws.Range("I102").Value = 1068.7778
$ws.Range("J102").Value = 3298
$ws.Range("K102").Value = 1068.7778
$ws.Range("L102").Value = 3298
$ws.Range("M102").Value = 553.2221999999999
$ws.Range("N102").Value = -6542
$ws.Range("H136").Value = 27521.076
$ws.Range("J136").Value = 27521.076
$ws.Range("L136").Value = 82563.228
$ws.Range("N136").Value = -87663.228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 4000
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H7").Value = 5397.6
$ws.Range("J7").Value = 6499.5
$ws.Range("L7").Value = 6499.5
$ws.Range("N7").Value = -6723.5
$ws.Range("H22").Value = 574.5
$ws.Range("I22").Value = 550
$ws.Range("K22").Value = 550
$ws.Range("M22").Value = -255
$ws.Range("H27").Value = 574.5
$ws.Range("I27").Value = 550
$ws.Range("K27").Value = 550
$ws.Range("M27").Value = -443
$ws.Range("H28").Value = 4000
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H37").Value = 4000
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H61").Value = 2625.375
$ws.Range("I61").Value = 3000.6
$ws.Range("K61").Value = 3000.6
$ws.Range("M61").Value = -2798.6
$ws.Range("H113").Value = 2625.375
$ws.Range("I113").Value = 3000.6
$ws.Range("K113").Value = 3000.6
$ws.Range("M113").Value = -830.5999999999999
$ws.Range("H122").Value = 6399.5
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 5397.6
$ws.Range("J126").Value = 6499.5
$ws.Range("L126").Value = 19498.5
$ws.Range("N126").Value = -24438.5
$ws.Range("H132").Value = 1999.6666
$ws.Range("I132").Value = 1499.5
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 4498.5
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -1968.5
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H126").Value = 1189.7778
$ws.Range("I126").Value = 1224.75
$ws.Range("K126").Value = 3674.25
$ws.Range("M126").Value = -1204.25
$ws.Range("H132").Value = 251211.5
$ws.Range("I132").Value = 251211.5
$ws.Range("K132").Value = 753634.5
$ws.Range("M132").Value = -751104.5
$ws.Range("H133").Value = 49997.332
$ws.Range("J133").Value = 49997.332
$ws.Range("L133").Value = 49997.332
$ws.Range("N133").Value = -60117.332
